# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" for the file that was just
# (re-)handed off (0cc65571-8b6f-48c1-ae04-253a6f8106cc), in each
# per-locale status sheet.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Re-stamp the handoff file name (unchanged text, but re-written as part
# of regenerating the handoff report) and the handoff datetime (new,
# later timestamp) for row 4 - the 0cc65571... entry.
$zhcn.Range("C4").Value = "0cc65571-8b6f-48c1-ae04-253a6f8106cc.bdf086cab8953438a9e19a960ed716252eed8fe8.zh-cn.xlf"
$zhcn.Range("D4").Value = "2016-03-04 02:13:06"

$dede.Range("C4").Value = "0cc65571-8b6f-48c1-ae04-253a6f8106cc.bdf086cab8953438a9e19a960ed716252eed8fe8.de-de.xlf"
$dede.Range("D4").Value = "2016-03-04 02:13:20"
